# Apply updated crypto price/volume data to the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.314.15"
$ws.Range("E2").Value = "  -0.16%  "
$ws.Range("D3").Value = "3.676.02"
$ws.Range("E3").Value = "  -0.43%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D5").Value = "'682.08"
$ws.Range("E5").Value = "  -0.39%  "
$ws.Range("D6").Value = "'158.52"
$ws.Range("E6").Value = "  -2.70%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -1.38%  "
$ws.Range("E9").Value = "  -1.75%  "
$ws.Range("E10").Value = "  -4.22%  "
$ws.Range("D11").Value = "'0.436"
$ws.Range("E11").Value = "  -3.30%  "
$ws.Range("E12").Value = "  -2.15%  "
$ws.Range("D13").Value = "4.296.44"
$ws.Range("E13").Value = "  -0.46%  "
$ws.Range("E14").Value = "  -4.19%  "
$ws.Range("D15").Value = "3.670.48"
$ws.Range("E15").Value = "  -0.45%  "
$ws.Range("D16").Value = "69.296.18"
$ws.Range("E16").Value = "  -0.27%  "
$ws.Range("E17").Value = "  +1.90%  "
$ws.Range("D18").Value = "'15.87"
$ws.Range("E18").Value = "  -2.86%  "
$ws.Range("E19").Value = "  -4.06%  "
$ws.Range("D20").Value = "'469.12"
$ws.Range("E20").Value = "  -2.39%  "
$ws.Range("D21").Value = "'9.97"
$ws.Range("E21").Value = "  +1.69%  "
$ws.Range("D22").Value = "'0.649"
$ws.Range("E22").Value = "  -2.67%  "
$ws.Range("D23").Value = "'79.97"
$ws.Range("E23").Value = "  -0.36%  "
$ws.Range("D24").Value = "3.821.93"
$ws.Range("E24").Value = "  -0.39%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("E26").Value = "  -6.05%  "
$ws.Range("D27").Value = "'10.90"
$ws.Range("E27").Value = "  -5.11%  "
$ws.Range("D28").Value = "'9.13"
$ws.Range("E28").Value = "  -4.94%  "
$ws.Range("E29").Value = "  -2.16%  "
$ws.Range("E30").Value = "  -5.43%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").Value = "'6.56"
$ws.Range("E31").Value = "  -4.39%  "
$ws.Range("B32").Value = "Binance-PegBSC-USD"
$ws.Range("C32").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("E33").Value = "  -5.92%  "
$ws.Range("D34").Value = "'26.80"
$ws.Range("E34").Value = "  -1.14%  "
$ws.Range("D35").Value = "3.654.65"
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("D36").Value = "'0.158"
$ws.Range("E36").Value = "  -3.78%  "
$ws.Range("D37").Value = "'8.16"
$ws.Range("E37").Value = "  -5.10%  "
$ws.Range("E38").Value = "  -1.07%  "
$ws.Range("E40").Value = "  +2.32%  "
$ws.Range("D41").Value = "'0.0900"
$ws.Range("E41").Value = "  -4.53%  "
$ws.Range("E42").Value = "  -0.08%  "
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").Value = "'166.45"
$ws.Range("E43").Value = "  +7.15%  "
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").Value = "'0.940"
$ws.Range("E44").Value = "  -2.24%  "
$ws.Range("D45").Value = "'47.60"
$ws.Range("E45").Value = "  -1.22%  "
$ws.Range("D46").Value = "'2.73"
$ws.Range("E46").Value = "  -4.21%  "
$ws.Range("E47").Value = "  -1.63%  "
$ws.Range("E48").Value = "  +0.63%  "
$ws.Range("E49").Value = "  -3.75%  "
$ws.Range("E50").Value = "  -4.48%  "
$ws.Range("D51").Value = "'26.97"
$ws.Range("E51").Value = "  -3.35%  "
